# Refresh the cryptocurrency price/volume snapshot (and fix the Hedera /
# TrustWalletToken row order, which swapped places in the source feed).
#
# Column D ("Price") holds literal text in this sheet: some values are not
# even valid numbers (e.g. "29.377.87", "7.900"), and the ones that do look
# numeric (e.g. "1.019") must still keep their exact original decimal text.
# Excel's normal Value-setter auto-coerces numeric-looking strings into
# binary floating point ("1.019" becomes 1.0189999999999999 once it round-
# trips through a double), so every Price cell is switched to the Text
# number format ("@") right before its literal string is written, which
# guarantees an exact, lossless assignment. Column E (percent change) and
# the Coin/Link columns are unambiguous text already and need no such
# treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.178.57' },
    @{ Cell = 'E2'; Value = '  -3.99%  ' },
    @{ Cell = 'D3'; Value = '1.967.68' },
    @{ Cell = 'E3'; Value = '  -6.28%  ' },
    @{ Cell = 'D4'; Value = '1.022' },
    @{ Cell = 'E4'; Value = '  +2.04%  ' },
    @{ Cell = 'D5'; Value = '327.23' },
    @{ Cell = 'E5'; Value = '  -4.75%  ' },
    @{ Cell = 'D6'; Value = '1.019' },
    @{ Cell = 'E6'; Value = '  +1.73%  ' },
    @{ Cell = 'D7'; Value = '0.4978' },
    @{ Cell = 'E7'; Value = '  -7.11%  ' },
    @{ Cell = 'D8'; Value = '0.4177' },
    @{ Cell = 'E8'; Value = '  -5.91%  ' },
    @{ Cell = 'D9'; Value = '53.49' },
    @{ Cell = 'E9'; Value = '  -2.44%  ' },
    @{ Cell = 'D10'; Value = '0.08834' },
    @{ Cell = 'E10'; Value = '  -6.11%  ' },
    @{ Cell = 'D11'; Value = '1.092' },
    @{ Cell = 'E11'; Value = '  -6.80%  ' },
    @{ Cell = 'D12'; Value = '2.097.55' },
    @{ Cell = 'E12'; Value = '  -5.13%  ' },
    @{ Cell = 'D13'; Value = '22.84' },
    @{ Cell = 'E13'; Value = '  -7.78%  ' },
    @{ Cell = 'D14'; Value = '7.842' },
    @{ Cell = 'E14'; Value = '  -8.45%  ' },
    @{ Cell = 'D15'; Value = '6.367' },
    @{ Cell = 'E15'; Value = '  -7.98%  ' },
    @{ Cell = 'D16'; Value = '1.022' },
    @{ Cell = 'E16'; Value = '  +1.96%  ' },
    @{ Cell = 'D17'; Value = '91.43' },
    @{ Cell = 'E17'; Value = '  -10.40%  ' },
    @{ Cell = 'D18'; Value = '0.00001094' },
    @{ Cell = 'E18'; Value = '  -5.84%  ' },
    @{ Cell = 'D19'; Value = '0.06732' },
    @{ Cell = 'E19'; Value = '  +0.62%  ' },
    @{ Cell = 'D20'; Value = '19.14' },
    @{ Cell = 'E20'; Value = '  -9.73%  ' },
    @{ Cell = 'D21'; Value = '1.018' },
    @{ Cell = 'E21'; Value = '  +1.52%  ' },
    @{ Cell = 'D22'; Value = '5.914' },
    @{ Cell = 'E22'; Value = '  -6.61%  ' },
    @{ Cell = 'D23'; Value = '29.292.61' },
    @{ Cell = 'E23'; Value = '  -3.77%  ' },
    @{ Cell = 'D24'; Value = '11.82' },
    @{ Cell = 'E24'; Value = '  -5.78%  ' },
    @{ Cell = 'D25'; Value = '2.309' },
    @{ Cell = 'E25'; Value = '  -0.48%  ' },
    @{ Cell = 'D26'; Value = '20.51' },
    @{ Cell = 'E26'; Value = '  -6.45%  ' },
    @{ Cell = 'D27'; Value = '155.29' },
    @{ Cell = 'E27'; Value = '  -4.68%  ' },
    @{ Cell = 'D28'; Value = '6.143' },
    @{ Cell = 'E28'; Value = '  -9.48%  ' },
    @{ Cell = 'D29'; Value = '2.267' },
    @{ Cell = 'E29'; Value = '  -10.23%  ' },
    @{ Cell = 'D30'; Value = '125.88' },
    @{ Cell = 'E30'; Value = '  -5.83%  ' },
    @{ Cell = 'D31'; Value = '1.037' },
    @{ Cell = 'E31'; Value = '  -9.44%  ' },
    @{ Cell = 'D32'; Value = '0.09845' },
    @{ Cell = 'E32'; Value = '  -6.88%  ' },
    @{ Cell = 'D33'; Value = '1.505' },
    @{ Cell = 'E33'; Value = '  -9.77%  ' },
    @{ Cell = 'D34'; Value = '5.742' },
    @{ Cell = 'E34'; Value = '  -8.35%  ' },
    @{ Cell = 'D35'; Value = '3.726' },
    @{ Cell = 'E35'; Value = '  -3.25%  ' },
    @{ Cell = 'D36'; Value = '0.02416' },
    @{ Cell = 'E36'; Value = '  -8.42%  ' },
    @{ Cell = 'D37'; Value = '9.116' },
    @{ Cell = 'E37'; Value = '  -10.54%  ' },
    @{ Cell = 'B38'; Value = 'TrustWalletToken' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D38'; Value = '1.285' },
    @{ Cell = 'E38'; Value = '  -4.76%  ' },
    @{ Cell = 'B39'; Value = 'Hedera' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D39'; Value = '0.06294' },
    @{ Cell = 'E39'; Value = '  -7.56%  ' },
    @{ Cell = 'D40'; Value = '0.6445' },
    @{ Cell = 'E40'; Value = '  -8.31%  ' },
    @{ Cell = 'D41'; Value = '11.44' },
    @{ Cell = 'E41'; Value = '  -9.73%  ' },
    @{ Cell = 'D42'; Value = '0.2009' },
    @{ Cell = 'E42'; Value = '  -9.60%  ' },
    @{ Cell = 'D43'; Value = '1.018' },
    @{ Cell = 'E43'; Value = '  +1.69%  ' },
    @{ Cell = 'D44'; Value = '0.6207' },
    @{ Cell = 'E44'; Value = '  -9.62%  ' },
    @{ Cell = 'D45'; Value = '13.44' },
    @{ Cell = 'E45'; Value = '  -7.09%  ' },
    @{ Cell = 'D46'; Value = '2.169' },
    @{ Cell = 'E46'; Value = '  -7.42%  ' },
    @{ Cell = 'D47'; Value = '1.273' },
    @{ Cell = 'E47'; Value = '  -8.94%  ' },
    @{ Cell = 'D48'; Value = '3.492' },
    @{ Cell = 'E48'; Value = '  -4.04%  ' },
    @{ Cell = 'D49'; Value = '0.00000000339' },
    @{ Cell = 'E49'; Value = '  -2.07%  ' },
    @{ Cell = 'D50'; Value = '0.06894' },
    @{ Cell = 'E50'; Value = '  -4.84%  ' },
    @{ Cell = 'D51'; Value = '1.109' },
    @{ Cell = 'E51'; Value = '  -9.91%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell -match '^D\d+$') {
        # Force Text format so the literal decimal string round-trips
        # byte-for-byte instead of being parsed into a Double.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
